# Value of a Statistical Life workbook update:
# - EPA's mortality-risk-valuation page moved; replace the old (dead) source
#   URL with the current one and turn it into a real hyperlink.
# - Enable iterative calculation (matches workbook calc options).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Turn on iterative calculation for the workbook.
$excel.Iteration = $true
$excel.MaxChange = 0.00001

$c = $ws.Range("B6")

# Add the new hyperlink (this also updates the cell's displayed text to the
# "TextToDisplay" argument, so we fix the cell text up afterwards).
$ws.Hyperlinks.Add(
    $c,
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation",
    "whatvalue",
    [Type]::Missing,
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation - whatvalue"
) | Out-Null

# Cell should display the full URL (with fragment) as its text.
$c.Value = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"

# Keep using the existing "Hyperlink" cell style rather than creating a
# duplicate style entry.
$c.Style = "Hyperlink"
